# Update column G ("K") values on Sheet1 per the diff:
# G2: 2 -> 0
# G3: 8 -> 2
# G4: 7 -> 1
# G5: 5 -> 0
# G6: 3 -> 1
# G7: 6 -> 0
# G8: 2 -> 1
# G10: 3 -> 1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 0
$ws.Range("G6").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 1
$ws.Range("G10").Value = 1
